$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the matched-error values in place (diagonal re-alignment caused by
# adding the full ifoCAST series into the evaluation window).
$ws.Range("B2").Value = 0.08873584232078371
$ws.Range("C2").Value = 0.5038494195916882
$ws.Range("D2").Value = -0.08040038373236996
$ws.Range("E2").Value = 0.7915874115568508
$ws.Range("F2").Value = 0.7713569177313746
$ws.Range("G2").Value = 0.3749198787210216
$ws.Range("H2").Value = 0.5168873173527261
$ws.Range("I2").Value = 0.6892627280777406
$ws.Range("J2").Value = 0.1844308218533179
$ws.Range("K2").Value = 0.4621062544735015

$ws.Range("B3").Value = 0.5262773992246967
$ws.Range("C3").Value = -0.07779684416992128
$ws.Range("D3").Value = 0.7695197259028355
$ws.Range("E3").Value = 0.7653749374270196
$ws.Range("F3").Value = 0.3695766285386105
$ws.Range("G3").Value = 0.5065187839823598
$ws.Range("H3").Value = 0.6807854059541167
$ws.Range("I3").Value = 0.1766585135465071
$ws.Range("J3").Value = 0.4534822809405329
$ws.Range("K3").Value = 0.2839429537672726

$ws.Range("B4").Value = -0.05637227271270118
$ws.Range("C4").Value = 0.8896975747020479
$ws.Range("D4").Value = 0.6761805930132365
$ws.Range("E4").Value = 0.3398496246900327
$ws.Range("F4").Value = 0.51514490401227
$ws.Range("G4").Value = 0.6548452325286815
$ws.Range("H4").Value = 0.1538799505182796
$ws.Range("I4").Value = 0.4398642457919151
$ws.Range("J4").Value = 0.2654221565820968
$ws.Range("K4").Value = 0.5706702220727796

$ws.Range("B5").Value = 0.8473769423816584
$ws.Range("C5").Value = 0.6469761191575033
$ws.Range("D5").Value = 0.3368290248851115
$ws.Range("E5").Value = 0.4962836204007895
$ws.Range("F5").Value = 0.6345141014634773
$ws.Range("G5").Value = 0.1393497950840318
$ws.Range("H5").Value = 0.4232232002996471
$ws.Range("I5").Value = 0.2478383111533466
$ws.Range("J5").Value = 0.5541963385427369
$ws.Range("K5").Value = -0.05728305666909728

$ws.Range("B6").Value = 0.9871738597754134
$ws.Range("C6").Value = 0.4122000866690486
$ws.Range("D6").Value = 0.3051903340877904
$ws.Range("E6").Value = 0.6592080140502106
$ws.Range("F6").Value = 0.1487291986403278
$ws.Range("G6").Value = 0.3684124016882561
$ws.Range("H6").Value = 0.2289442202646768
$ws.Range("I6").Value = 0.5393323377276911
$ws.Range("J6").Value = -0.085994767298468
$ws.Range("K6").Value = 0.6071338394308724

$ws.Range("B7").Value = 0.8628949586592991
$ws.Range("C7").Value = 0.3523010363001488
$ws.Range("D7").Value = 0.4189247832594023
$ws.Range("E7").Value = 0.1846742797061906
$ws.Range("F7").Value = 0.3663616442486632
$ws.Range("G7").Value = 0.1523759808286466
$ws.Range("H7").Value = 0.5119329433524077
$ws.Range("I7").Value = -0.1122838174961934
$ws.Range("J7").Value = 0.5637365487175399
$ws.Range("K7").ClearContents()

$ws.Range("B8").Value = 0.6646266232236873
$ws.Range("C8").Value = 0.5522135229949265
$ws.Range("D8").Value = 0.005595857889999001
$ws.Range("E8").Value = 0.3945558830192304
$ws.Range("F8").Value = 0.1884012920210401
$ws.Range("G8").Value = 0.4743913731481941
$ws.Range("H8").Value = -0.1178261565835546
$ws.Range("I8").Value = 0.5679118776562884
$ws.Range("J8").ClearContents()

$ws.Range("B9").Value = 0.7878040141027678
$ws.Range("C9").Value = 0.09027459876430857
$ws.Range("D9").Value = 0.2489555163855133
$ws.Range("E9").Value = 0.1983550348802827
$ws.Range("F9").Value = 0.489756542847739
$ws.Range("G9").Value = -0.1566744915409777
$ws.Range("H9").Value = 0.5520660686291026
$ws.Range("I9").ClearContents()

$ws.Range("B10").Value = 0.4012987852456914
$ws.Range("C10").Value = 0.3660442907967085
$ws.Range("D10").Value = 0.03589089343071589
$ws.Range("E10").Value = 0.5187154933129405
$ws.Range("F10").Value = -0.1209316190860741
$ws.Range("G10").Value = 0.5207384222454754
$ws.Range("H10").ClearContents()

$ws.Range("B11").Value = 0.6128694682008229
$ws.Range("C11").Value = 0.05354676096860539
$ws.Range("D11").Value = 0.4240932542019461
$ws.Range("E11").Value = -0.088714230837594
$ws.Range("F11").Value = 0.5331265480731927
$ws.Range("G11").ClearContents()

$ws.Range("B12").Value = 0.2932231203848173
$ws.Range("C12").Value = 0.5091174976711597
$ws.Range("D12").Value = -0.2051988091956081
$ws.Range("E12").Value = 0.5482318197250452
$ws.Range("F12").ClearContents()

$ws.Range("B13").Value = 0.6739775747052469
$ws.Range("C13").Value = -0.1916615369476919
$ws.Range("D13").Value = 0.4852589007350822
$ws.Range("E13").ClearContents()

$ws.Range("B14").Value = 0.06218750491771613
$ws.Range("C14").Value = 0.5845769509171186
$ws.Range("D14").ClearContents()

$ws.Range("B15").Value = 0.6286366421565677
$ws.Range("C15").ClearContents()

$ws.Range("B16").ClearContents()
